# The task "only allow one instance of this program open at a time" (Id 15)
# is moved from the "Active" sheet (status Todo) to the "Inactive" sheet,
# marking it as "Skip" with a Done date, as if the user re-filed the row
# between the two tracking tabs.

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Capture the row's data before removing it from "Active".
$taskId = $active.Range("A4").Value()
$title = $active.Range("B4").Value()
$category = $active.Range("D4").Value()
$created = $active.Range("E4").Value()

# Remove the row from "Active" - remaining rows shift up.
$active.Rows(4).Delete()

# Insert a new row in "Inactive" right after the header, pushing the
# existing rows down.
$inactive.Rows(2).Insert()
$inactive.Range("A2:F2").ClearFormats()

# Format the date-looking columns as text first so Excel doesn't silently
# convert "12/15/2017" / "3/5/2018" into date serials - every other row in
# this column stores them as plain text.
$inactive.Range("E2:F2").NumberFormat = "@"

$inactive.Range("A2").Value = $taskId
$inactive.Range("B2").Value = $title
$inactive.Range("C2").Value = "Skip"
$inactive.Range("D2").Value = $category
$inactive.Range("E2").Value = $created
$inactive.Range("F2").Value = "3/5/2018"

# Drop the text-format override so the new cells match the plain default
# styling (s="0") used by the rest of the data rows.
$inactive.Range("A2:F2").ClearFormats()
